$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new Price values are plain decimal numbers (e.g. "0.579",
# "6.70", "0.0991"). Mark just those specific cells as Text first so Excel
# stores the digits verbatim -- matching the source workbook, where every
# Price/Volume cell is a string, not a number -- instead of silently
# converting/rounding them into binary-float numeric cells.
$textCells = @('D4', 'D5', 'D6', 'D7', 'D8', 'D10', 'D19', 'D20', 'D21', 'D23', 'D24', 'D27', 'D29', 'D30', 'D33', 'D34', 'D35', 'D36', 'D37', 'D40', 'D42', 'D43', 'D44', 'D45', 'D46', 'D48')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '60.907.42'
$ws.Range('E2').Value = '  +3.13%  '
$ws.Range('D3').Value = '2.703.20'
$ws.Range('E3').Value = '  +2.19%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').Value = '526.42'
$ws.Range('E5').Value = '  +1.47%  '
$ws.Range('D6').Value = '145.02'
$ws.Range('E6').Value = '  -0.43%  '
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  -0.28%  '
$ws.Range('D8').Value = '0.579'
$ws.Range('E8').Value = '  +2.22%  '
$ws.Range('D9').Value = '2.727.90'
$ws.Range('E9').Value = '  +2.54%  '
$ws.Range('D10').Value = '6.70'
$ws.Range('E10').Value = '  +6.93%  '
$ws.Range('E11').Value = '  +1.00%  '
$ws.Range('E12').Value = '  +1.08%  '
$ws.Range('E13').Value = '  +2.96%  '
$ws.Range('D14').Value = '3.181.61'
$ws.Range('E14').Value = '  +1.86%  '
$ws.Range('D15').Value = '60.854.42'
$ws.Range('E15').Value = '  +2.80%  '
$ws.Range('D16').Value = '2.892.34'
$ws.Range('E16').Value = '  +8.84%  '
$ws.Range('E17').Value = '  +1.67%  '
$ws.Range('E18').Value = '  +0.59%  '
$ws.Range('D19').Value = '348.12'
$ws.Range('E19').Value = '  -0.53%  '
$ws.Range('D20').Value = '4.51'
$ws.Range('E20').Value = '  -0.26%  '
$ws.Range('D21').Value = '10.61'
$ws.Range('E21').Value = '  +2.42%  '
$ws.Range('E22').Value = '  +3.57%  '
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').Value = '63.82'
$ws.Range('E24').Value = '  +2.63%  '
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('E26').Value = '  +4.70%  '
$ws.Range('D27').Value = '0.993'
$ws.Range('E27').Value = '  -0.38%  '
$ws.Range('E28').Value = '  +1.48%  '
$ws.Range('D29').Value = '7.30'
$ws.Range('E29').Value = '  +2.20%  '
$ws.Range('D30').Value = '6.76'
$ws.Range('E30').Value = '  +7.89%  '
$ws.Range('E31').Value = '  -0.22%  '
$ws.Range('E32').Value = '  +1.66%  '
$ws.Range('D33').Value = '19.15'
$ws.Range('E33').Value = '  +0.70%  '
$ws.Range('D34').Value = '150.26'
$ws.Range('E34').Value = '  +0.18%  '
$ws.Range('D35').Value = '4.24'
$ws.Range('E35').Value = '  +5.38%  '
$ws.Range('B36').Value = 'SuiNetwork'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D36').Value = '0.948'
$ws.Range('E36').Value = '  -1.47%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '1.23'
$ws.Range('E37').Value = '  +7.77%  '
$ws.Range('E38').Value = '  +4.17%  '
$ws.Range('E39').Value = '  +7.50%  '
$ws.Range('D40').Value = '36.99'
$ws.Range('E40').Value = '  +0.66%  '
$ws.Range('E41').Value = '  -1.02%  '
$ws.Range('D42').Value = '285.03'
$ws.Range('E42').Value = '  +3.00%  '
$ws.Range('D43').Value = '20.19'
$ws.Range('E43').Value = '  +2.79%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').Value = '0.0991'
$ws.Range('E44').Value = '  +0.70%  '
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').Value = '0.612'
$ws.Range('E45').Value = '  +0.53%  '
$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D46').Value = '0.997'
$ws.Range('E46').Value = '  +0.12%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '2.139.69'
$ws.Range('E47').Value = '  +7.30%  '
$ws.Range('D48').Value = '0.0539'
$ws.Range('E48').Value = '  +3.21%  '
$ws.Range('E49').Value = '  +2.16%  '
$ws.Range('E50').Value = '  +2.95%  '
$ws.Range('E51').Value = '  +1.72%  '
